# Insert a new weekly record above the current row 51 ("Fruta / hortaliza,
# semanal"). This pushes the existing rows 51-66 down to 52-67 and fills
# the freshly inserted row 51 with the new Arica y Parinacota observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 51..66 down to 52..67, leaving row 51 blank (but formatted,
# since Insert() carries the formatting of the row above/below as Excel does).
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with the new data point.
$ws.Cells.Item(51, 1).Value  = 4
$ws.Cells.Item(51, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(51, 3).Value  = "Los Lagos"
$ws.Cells.Item(51, 4).Value  = 44463
$ws.Cells.Item(51, 5).Value  = 10
$ws.Cells.Item(51, 6).Value  = 100112052
$ws.Cells.Item(51, 7).Value  = "Albahaca"
$ws.Cells.Item(51, 8).Value  = "Sin especificar"
$ws.Cells.Item(51, 9).Value  = "Primera"
$ws.Cells.Item(51, 10).Value = 90
$ws.Cells.Item(51, 11).Value = 6000
$ws.Cells.Item(51, 12).Value = 6000
$ws.Cells.Item(51, 13).Value = 6000
$ws.Cells.Item(51, 14).Value = "$/paquete"
$ws.Cells.Item(51, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(51, 16).Value = 6000
$ws.Cells.Item(51, 17).Value = 1
$ws.Cells.Item(51, 18).Value = "Hortaliza"
